# Re-pull / push updated data: update the dSF (column F) values for several rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0
$ws.Range("F11").Value = -2
$ws.Range("F14").Value = -4
$ws.Range("F15").Value = -1
$ws.Range("F19").Value = -2
